$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.316.58'
$ws.Range("E2").Value = '  +2.18%  '

$ws.Range("D3").Value = '2.500.59'
$ws.Range("E3").Value = '  +1.81%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''323.46'
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("D6").Value = '''109.08'
$ws.Range("E6").Value = '  +3.47%  '

$ws.Range("E7").Value = '  +1.18%  '

$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '''0.536'
$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").Value = '''39.15'
$ws.Range("E10").Value = '  +8.58%  '

$ws.Range("E11").Value = '  +0.63%  '

$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").Value = '''18.38'
$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").Value = '''7.19'
$ws.Range("E14").Value = '  +1.37%  '

$ws.Range("D15").Value = '2.890.13'
$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("D16").Value = '2.506.01'
$ws.Range("E16").Value = '  +1.36%  '

$ws.Range("E17").Value = '  +1.20%  '

$ws.Range("D18").Value = '47.239.19'
$ws.Range("E18").Value = '  +2.37%  '

$ws.Range("D19").Value = '''12.85'
$ws.Range("E19").Value = '  +1.70%  '

$ws.Range("D20").Value = '''6.63'
$ws.Range("E20").Value = '  +3.20%  '

$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("E22").Value = '  +12.76%  '

$ws.Range("D23").Value = '''70.75'
$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("D24").Value = '''247.88'
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("E25").Value = '  +3.11%  '

$ws.Range("D26").Value = '''26.09'
$ws.Range("E26").Value = '  +0.38%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("E28").Value = '  +4.42%  '

$ws.Range("D29").Value = '''10.07'
$ws.Range("E29").Value = '  +3.73%  '

$ws.Range("D30").Value = '''35.25'
$ws.Range("E30").Value = '  +2.35%  '

$ws.Range("E31").Value = '  +7.06%  '

$ws.Range("D32").Value = '''49.81'
$ws.Range("E32").Value = '  +0.91%  '

$ws.Range("D33").Value = '''20.07'
$ws.Range("E33").Value = '  +1.04%  '

$ws.Range("E34").Value = '  +1.65%  '

$ws.Range("E35").Value = '  +3.04%  '

$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("E37").Value = '  +3.82%  '

$ws.Range("E38").Value = '  +2.95%  '

$ws.Range("E39").Value = '  +1.03%  '

$ws.Range("E40").Value = '  +1.24%  '

$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''121.94'
$ws.Range("E41").Value = '  -3.46%  '

$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '''2.24'
$ws.Range("E42").Value = '  +0.40%  '

$ws.Range("E43").Value = '  +2.10%  '

$ws.Range("E44").Value = '  +2.16%  '

$ws.Range("D45").Value = '1.991.52'
$ws.Range("E45").Value = '  +0.93%  '

$ws.Range("E46").Value = '  +2.63%  '

$ws.Range("E47").Value = '  -1.10%  '

$ws.Range("D48").Value = '''1.78'
$ws.Range("E48").Value = '  -4.42%  '

$ws.Range("E49").Value = '  -0.34%  '

$ws.Range("E50").Value = '  +3.05%  '

$ws.Range("D51").Value = '''56.55'
$ws.Range("E51").Value = '  +3.73%  '
